$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '38.764.49'
$ws.Range("D2").NumberFormat = "General"
$ws.Range("E2").Value = '  +1.91%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.094.09'
$ws.Range("D3").NumberFormat = "General"
$ws.Range("E3").Value = '  +0.20%  '
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '229.77'
$ws.Range("D5").NumberFormat = "General"
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("E6").Value = '  +0.40%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '61.25'
$ws.Range("D7").NumberFormat = "General"
$ws.Range("E7").Value = '  +0.45%  '
$ws.Range("E8").Value = '  -0.02%  '
$ws.Range("E9").Value = '  +1.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0845'
$ws.Range("D10").NumberFormat = "General"
$ws.Range("E10").Value = '  +0.65%  '
$ws.Range("B11").Value = 'TRON'
$ws.Range("C11").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.104'
$ws.Range("D11").NumberFormat = "General"
$ws.Range("E11").Value = '  -0.09%  '
$ws.Range("B12").Value = 'Chainlink'
$ws.Range("C12").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.41'
$ws.Range("D12").NumberFormat = "General"
$ws.Range("E12").Value = '  +5.31%  '
$ws.Range("B13").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C13").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.446.92'
$ws.Range("D13").NumberFormat = "General"
$ws.Range("E13").Value = '  +1.99%  '
$ws.Range("E14").Value = '  -0.33%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.808'
$ws.Range("D15").NumberFormat = "General"
$ws.Range("E15").Value = '  +4.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.47'
$ws.Range("D16").NumberFormat = "General"
$ws.Range("E16").Value = '  -0.30%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.095.32'
$ws.Range("D17").NumberFormat = "General"
$ws.Range("E17").Value = '  -0.96%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '38.722.56'
$ws.Range("D18").NumberFormat = "General"
$ws.Range("E18").Value = '  +1.95%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '72.00'
$ws.Range("D19").NumberFormat = "General"
$ws.Range("E19").Value = '  +2.87%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.09'
$ws.Range("D20").NumberFormat = "General"
$ws.Range("E20").Value = '  +0.69%  '
$ws.Range("E21").Value = '  +0.45%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.88'
$ws.Range("D22").NumberFormat = "General"
$ws.Range("E22").Value = '  +1.73%  '
$ws.Range("E23").Value = '  -0.43%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.38'
$ws.Range("D24").NumberFormat = "General"
$ws.Range("E24").Value = '  -2.01%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.35'
$ws.Range("D25").NumberFormat = "General"
$ws.Range("E25").Value = '  +1.55%  '
$ws.Range("B26").Value = 'Monero'
$ws.Range("C26").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '171.31'
$ws.Range("D26").NumberFormat = "General"
$ws.Range("E26").Value = '  +1.04%  '
$ws.Range("B27").Value = 'Cosmos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.56'
$ws.Range("D27").NumberFormat = "General"
$ws.Range("E27").Value = '  +1.85%  '
$ws.Range("E28").Value = '  +5.15%  '
$ws.Range("E29").Value = '  +7.74%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.30'
$ws.Range("D30").NumberFormat = "General"
$ws.Range("E30").Value = '  +1.85%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.49'
$ws.Range("D31").NumberFormat = "General"
$ws.Range("E31").Value = '  +4.34%  '
$ws.Range("E32").Value = '  +0.62%  '
$ws.Range("E33").Value = '  +1.95%  '
$ws.Range("E34").Value = '  +1.11%  '
$ws.Range("E35").Value = '  +0.69%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.59'
$ws.Range("D36").NumberFormat = "General"
$ws.Range("E36").Value = '  +2.05%  '
$ws.Range("E37").Value = '  -0.50%  '
$ws.Range("E38").Value = '  +1.50%  '
$ws.Range("E39").Value = '  -0.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.04'
$ws.Range("D40").NumberFormat = "General"
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  +5.29%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '101.12'
$ws.Range("D42").NumberFormat = "General"
$ws.Range("E42").Value = '  +1.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.534.25'
$ws.Range("D43").NumberFormat = "General"
$ws.Range("E43").Value = '  -0.64%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.80'
$ws.Range("D44").NumberFormat = "General"
$ws.Range("E44").Value = '  -0.70%  '
$ws.Range("E45").Value = '  +1.08%  '
$ws.Range("E46").Value = '  +1.54%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.67'
$ws.Range("D47").NumberFormat = "General"
$ws.Range("E47").Value = '  +5.92%  '
$ws.Range("E48").Value = '  -1.06%  '
$ws.Range("E49").Value = '  +1.62%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.98'
$ws.Range("D50").NumberFormat = "General"
$ws.Range("E50").Value = '  -0.71%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.289.97'
$ws.Range("D51").NumberFormat = "General"
$ws.Range("E51").Value = '  +0.15%  '
